$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2318059299191375
$ws.Cells.Item(2, 3).Value = 0.5148247978436657
$ws.Cells.Item(2, 10).Value = 0.01347708894878706
$ws.Cells.Item(2, 15).Value = 0.002695417789757413
$ws.Cells.Item(2, 16).Value = 0.1536388140161725
$ws.Cells.Item(2, 19).Value = 0.08355795148247978
$ws.Cells.Item(3, 2).Value = 0.01470588235294118
$ws.Cells.Item(3, 3).Value = 0.02450980392156863
$ws.Cells.Item(3, 16).Value = 0.7058823529411765
$ws.Cells.Item(3, 19).Value = 0.2254901960784314
$ws.Cells.Item(4, 10).Value = 0.05882352941176471
$ws.Cells.Item(4, 16).Value = 0.7843137254901961
$ws.Cells.Item(4, 19).Value = 0.1568627450980392
$ws.Cells.Item(5, 16).Value = 0.5
$ws.Cells.Item(5, 19).Value = 0.5
$ws.Cells.Item(6, 2).Value = 0.104089219330855
$ws.Cells.Item(6, 4).Value = 0.01486988847583643
$ws.Cells.Item(6, 6).Value = 0.104089219330855
$ws.Cells.Item(6, 10).Value = 0.20817843866171
$ws.Cells.Item(6, 15).Value = 0.04460966542750929
$ws.Cells.Item(6, 17).Value = 0.1263940520446097
$ws.Cells.Item(6, 18).Value = 0.09665427509293681
$ws.Cells.Item(6, 19).Value = 0.3011152416356878
$ws.Cells.Item(7, 2).Value = 0.1420454545454546
$ws.Cells.Item(7, 4).Value = 0.01136363636363636
$ws.Cells.Item(7, 6).Value = 0.08522727272727272
$ws.Cells.Item(7, 10).Value = 0.1079545454545455
$ws.Cells.Item(7, 15).Value = 0.02272727272727273
$ws.Cells.Item(7, 17).Value = 0.125
$ws.Cells.Item(7, 18).Value = 0.09659090909090909
$ws.Cells.Item(7, 19).Value = 0.4090909090909091
$ws.Cells.Item(8, 2).Value = 0.1157635467980296
$ws.Cells.Item(8, 4).Value = 0.01970443349753695
$ws.Cells.Item(8, 6).Value = 0.05911330049261083
$ws.Cells.Item(8, 10).Value = 0.125615763546798
$ws.Cells.Item(8, 15).Value = 0.01477832512315271
$ws.Cells.Item(8, 17).Value = 0.1428571428571428
$ws.Cells.Item(8, 18).Value = 0.07881773399014778
$ws.Cells.Item(8, 19).Value = 0.4433497536945813
$ws.Cells.Item(9, 2).Value = 0.124223602484472
$ws.Cells.Item(9, 4).Value = 0.0124223602484472
$ws.Cells.Item(9, 5).Value = 0.006211180124223602
$ws.Cells.Item(9, 6).Value = 0.08695652173913043
$ws.Cells.Item(9, 10).Value = 0.06211180124223602
$ws.Cells.Item(9, 15).Value = 0.02484472049689441
$ws.Cells.Item(9, 17).Value = 0.1490683229813665
$ws.Cells.Item(9, 18).Value = 0.1055900621118012
$ws.Cells.Item(9, 19).Value = 0.4285714285714285
$ws.Cells.Item(10, 2).Value = 0.1300081766148814
$ws.Cells.Item(10, 4).Value = 0.02861815208503679
$ws.Cells.Item(10, 5).Value = 0.001635322976287817
$ws.Cells.Item(10, 6).Value = 0.0776778413736713
$ws.Cells.Item(10, 10).Value = 0.1210139002452984
$ws.Cells.Item(10, 15).Value = 0.02289452166802943
$ws.Cells.Item(10, 17).Value = 0.169255928045789
$ws.Cells.Item(10, 18).Value = 0.08503679476696648
$ws.Cells.Item(10, 19).Value = 0.3638593622240393
$ws.Cells.Item(11, 7).Value = 0.1326530612244898
$ws.Cells.Item(11, 10).Value = 0.1020408163265306
$ws.Cells.Item(11, 11).Value = 0.2244897959183673
$ws.Cells.Item(11, 12).Value = 0.5340136054421769
$ws.Cells.Item(11, 19).Value = 0.006802721088435374
$ws.Cells.Item(12, 7).Value = 0.7160493827160493
$ws.Cells.Item(12, 10).Value = 0.2160493827160494
$ws.Cells.Item(12, 11).Value = 0.006172839506172839
$ws.Cells.Item(12, 12).Value = 0.0308641975308642
$ws.Cells.Item(12, 19).Value = 0.0308641975308642
$ws.Cells.Item(13, 7).Value = 0.6585365853658537
$ws.Cells.Item(13, 10).Value = 0.3170731707317073
$ws.Cells.Item(13, 19).Value = 0.02439024390243903
$ws.Cells.Item(15, 6).Value = 0.03211009174311927
$ws.Cells.Item(15, 8).Value = 0.1330275229357798
$ws.Cells.Item(15, 9).Value = 0.06880733944954129
$ws.Cells.Item(15, 10).Value = 0.3899082568807339
$ws.Cells.Item(15, 11).Value = 0.07798165137614679
$ws.Cells.Item(15, 13).Value = 0.01376146788990826
$ws.Cells.Item(15, 14).Value = 0.004587155963302753
$ws.Cells.Item(15, 15).Value = 0.06422018348623854
$ws.Cells.Item(15, 19).Value = 0.2155963302752294
$ws.Cells.Item(16, 6).Value = 0.02127659574468085
$ws.Cells.Item(16, 8).Value = 0.2212765957446808
$ws.Cells.Item(16, 9).Value = 0.08085106382978724
$ws.Cells.Item(16, 10).Value = 0.3404255319148936
$ws.Cells.Item(16, 11).Value = 0.1106382978723404
$ws.Cells.Item(16, 13).Value = 0.02553191489361702
$ws.Cells.Item(16, 15).Value = 0.0425531914893617
$ws.Cells.Item(16, 19).Value = 0.1574468085106383
$ws.Cells.Item(17, 6).Value = 0.02949852507374631
$ws.Cells.Item(17, 8).Value = 0.191740412979351
$ws.Cells.Item(17, 9).Value = 0.06784660766961652
$ws.Cells.Item(17, 10).Value = 0.4454277286135693
$ws.Cells.Item(17, 11).Value = 0.1061946902654867
$ws.Cells.Item(17, 13).Value = 0.01474926253687316
$ws.Cells.Item(17, 14).Value = 0.002949852507374631
$ws.Cells.Item(17, 15).Value = 0.05604719764011799
$ws.Cells.Item(17, 19).Value = 0.08554572271386431
$ws.Cells.Item(18, 6).Value = 0.04123711340206185
$ws.Cells.Item(18, 8).Value = 0.1907216494845361
$ws.Cells.Item(18, 9).Value = 0.05154639175257732
$ws.Cells.Item(18, 10).Value = 0.4329896907216495
$ws.Cells.Item(18, 11).Value = 0.07731958762886598
$ws.Cells.Item(18, 13).Value = 0.01030927835051546
$ws.Cells.Item(18, 15).Value = 0.06701030927835051
$ws.Cells.Item(18, 19).Value = 0.1288659793814433
$ws.Cells.Item(19, 6).Value = 0.02775441547518924
$ws.Cells.Item(19, 8).Value = 0.1925988225399495
$ws.Cells.Item(19, 9).Value = 0.07821698906644239
$ws.Cells.Item(19, 10).Value = 0.3877207737594617
$ws.Cells.Item(19, 11).Value = 0.1093355761143818
$ws.Cells.Item(19, 13).Value = 0.02523128679562658
$ws.Cells.Item(19, 15).Value = 0.07401177460050462
$ws.Cells.Item(19, 19).Value = 0.1051303616484441

Write-Host "Updated 112 cells"